$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 816.17645
$ws.Range("I19").Value = 919.9091
$ws.Range("K19").Value = 919.9091
$ws.Range("M19").Value = -744.9091
$ws.Range("H41").Value = 2015.6666
$ws.Range("I41").Value = 2158.8
$ws.Range("K41").Value = 2158.8
$ws.Range("M41").Value = -1718.8
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -40972
$ws.Range("H98").Value = 2276.7856
$ws.Range("I98").Value = 2270.4443
$ws.Range("J98").Value = 2288.2
$ws.Range("K98").Value = 2270.4443
$ws.Range("L98").Value = 2288.2
$ws.Range("M98").Value = -772.4443000000001
$ws.Range("N98").Value = -5284.2
$ws.Range("H111").Value = 2629
$ws.Range("I111").Value = 2629
$ws.Range("K111").Value = 7887
$ws.Range("M111").Value = -4820
$ws.Range("H116").Value = 6141.7144
$ws.Range("I116").Value = 6298.8
$ws.Range("J116").Value = 5749
$ws.Range("K116").Value = 6298.8
$ws.Range("L116").Value = 5749
$ws.Range("M116").Value = -2856.8
$ws.Range("N116").Value = -12633
$ws.Range("H122").Value = 2276.7856
$ws.Range("I122").Value = 2270.4443
$ws.Range("J122").Value = 2288.2
$ws.Range("K122").Value = 6811.3329
$ws.Range("L122").Value = 6864.599999999999
$ws.Range("M122").Value = -4361.3329
$ws.Range("N122").Value = -11764.6
$ws.Range("H125").Value = 601.75
$ws.Range("I125").Value = 572
$ws.Range("J125").Value = 631.5
$ws.Range("K125").Value = 5148
$ws.Range("L125").Value = 5683.5
$ws.Range("M125").Value = -2688
$ws.Range("N125").Value = -10603.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4247.8125
$ws.Range("I32").Value = 3197.6667
$ws.Range("K32").Value = 3197.6667
$ws.Range("M32").Value = -2910.6667
$ws.Range("H122").Value = 6940.2
$ws.Range("I122").Value = 6940.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 20820.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -18370.6
$ws.Range("N122").Value = ""
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2261.125
$ws.Range("I94").Value = 1667.6923
$ws.Range("K94").Value = 1667.6923
$ws.Range("M94").Value = -1216.6923
$ws.Range("H99").Value = 1771.125
$ws.Range("I99").Value = 1061.6666
$ws.Range("K99").Value = 1061.6666
$ws.Range("M99").Value = 436.3334
$ws.Range("H105").Value = 2875
$ws.Range("I105").Value = 2833.3333
$ws.Range("K105").Value = 2833.3333
$ws.Range("M105").Value = -1086.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12175
$ws.Range("I86").Value = 4566.6665
$ws.Range("K86").Value = 4566.6665
$ws.Range("M86").Value = -3443.6665
$ws.Range("H89").Value = 12175
$ws.Range("I89").Value = 4566.6665
$ws.Range("K89").Value = 22833.3325
$ws.Range("M89").Value = -17217.3325
$ws.Range("H99").Value = 6718.6665
$ws.Range("I99").Value = 6718.6665
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6718.6665
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5220.6665
$ws.Range("N99").Value = ""
$ws.Range("H122").Value = 1755.5
$ws.Range("I122").Value = 670
$ws.Range("J122").Value = 1910.5714
$ws.Range("K122").Value = 2010
$ws.Range("L122").Value = 5731.7142
$ws.Range("M122").Value = 440
$ws.Range("N122").Value = -10631.7142
$ws.Range("H126").Value = 6718.6665
$ws.Range("I126").Value = 6718.6665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20155.9995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17685.9995
$ws.Range("N126").Value = ""
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 163
$ws.Range("I38").Value = 43.75
$ws.Range("K38").Value = 131.25
$ws.Range("M38").Value = 215.75
$ws.Range("H97").Value = 2003.5
$ws.Range("J97").Value = 2003.5
$ws.Range("L97").Value = 6010.5
$ws.Range("N97").Value = -7002.5
$ws.Range("H107").Value = 205
$ws.Range("J107").Value = 210
$ws.Range("L107").Value = 630
$ws.Range("N107").Value = -4470
$ws.Range("H131").Value = 2714
$ws.Range("I131").Value = 1999.75
$ws.Range("J131").Value = 3666.3333
$ws.Range("K131").Value = 5999.25
$ws.Range("L131").Value = 10998.9999
$ws.Range("M131").Value = -959.25
$ws.Range("N131").Value = -21078.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 110005.336
$ws.Range("J36").Value = 14999.5
$ws.Range("L36").Value = 14999.5
$ws.Range("N36").Value = -15969.5
$ws.Range("H70").Value = 12932.667
$ws.Range("I70").Value = 12449.5
$ws.Range("J70").Value = 13899
$ws.Range("K70").Value = 12449.5
$ws.Range("L70").Value = 13899
$ws.Range("M70").Value = -12179.5
$ws.Range("N70").Value = -14439
$ws.Range("H73").Value = 12932.667
$ws.Range("I73").Value = 12449.5
$ws.Range("J73").Value = 13899
$ws.Range("K73").Value = 12449.5
$ws.Range("L73").Value = 13899
$ws.Range("M73").Value = -11513.5
$ws.Range("N73").Value = -15771
$ws.Range("H102").Value = 1889
$ws.Range("I102").Value = 1533.5
$ws.Range("J102").Value = 2600
$ws.Range("K102").Value = 1533.5
$ws.Range("L102").Value = 2600
$ws.Range("M102").Value = 88.5
$ws.Range("N102").Value = -5844
$ws.Range("H122").Value = 2013.8182
$ws.Range("I122").Value = 1990.6666
$ws.Range("K122").Value = 5971.9998
$ws.Range("M122").Value = -3521.9998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 500
$ws.Range("K2").Value = 500
$ws.Range("M2").Value = -388
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1093
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -50980
$ws.Range("H48").Value = 30833
$ws.Range("I48").Value = 31249.5
$ws.Range("K48").Value = 31249.5
$ws.Range("M48").Value = -30588.5
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50466
$ws.Range("H82").Value = 1569
$ws.Range("I82").Value = 624.75
$ws.Range("K82").Value = 624.75
$ws.Range("M82").Value = -263.75
$ws.Range("H85").Value = 1569
$ws.Range("I85").Value = 624.75
$ws.Range("K85").Value = 624.75
$ws.Range("M85").Value = 623.25
$ws.Range("H122").Value = 3199.6
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -15397
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 18000
$ws.Range("J10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("N10").Value = -18338
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H93").Value = 29999
$ws.Range("J93").Value = 29999
$ws.Range("L93").Value = 29999
$ws.Range("N93").Value = -34991
$ws.Range("H107").Value = 465.33334
$ws.Range("I107").Value = 561.3333
$ws.Range("K107").Value = 1683.9999
$ws.Range("M107").Value = 236.0001
$ws.Range("H113").Value = 419.33334
$ws.Range("I113").Value = 310.5
$ws.Range("K113").Value = 931.5
$ws.Range("M113").Value = 1238.5
$ws.Range("H126").Value = 4953.875
$ws.Range("I126").Value = 4791.3076
$ws.Range("J126").Value = 5658.3335
$ws.Range("K126").Value = 14373.9228
$ws.Range("L126").Value = 16975.0005
$ws.Range("M126").Value = -11903.9228
$ws.Range("N126").Value = -21915.0005
